$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column I ("UserDashboardAccountId") currently holds the text "NULL" for
# every data row (2-43). The edit replaces that placeholder text with the
# numeric value 1 across the whole column.
for ($r = 2; $r -le 43; $r++) {
    $ws.Cells.Item($r, 9).Value = 1
}

# Update the saved view/selection state of the sheet (scroll position +
# active cell / selection).
$win = $excel.ActiveWindow
$win.ScrollRow = 18
$win.ScrollColumn = 1
$ws.Range("I35").Select()
